# Calendar.xlsx edit: "Add function to detect conflicting courses."
#
# This extends the time axis at the bottom of the calendar with new
# 15-minute slots (5:15 - 7:30), removes the three old course blocks
# (GS-GS-6600, GS-GS-6400, GS-NE-6112) and adds two new course blocks:
#   - GS-CC-6208  10:00-11:00  N311   (columns B, D, F ; rows 10-14)
#   - GS-CC-6202  11:00-12:30  M616   (columns C, E      ; rows 14-20)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Extend the time axis in column A with new 15-minute slots down to
#    7:30 (rows 39-48).
# ---------------------------------------------------------------------
$newTimes = @("5:15","5:30","5:45","6:00","6:15","6:30","6:45","7:00","7:15","7:30")
for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = 39 + $i
    $ws.Cells.Item($row, 1).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------
# 2. Remove the old course blocks entirely (values, styles and merges).
# ---------------------------------------------------------------------
$ws.Range("C6:C18").UnMerge()
$ws.Range("C6:C18").Clear()

$ws.Range("B23:B27").UnMerge()
$ws.Range("D23:D27").UnMerge()
$ws.Range("F23:F27").UnMerge()
$ws.Range("B23:B27").Clear()
$ws.Range("D23:D27").Clear()
$ws.Range("F23:F27").Clear()

$ws.Range("B28:B32").UnMerge()
$ws.Range("E28:E32").UnMerge()
$ws.Range("B28:B32").Clear()
$ws.Range("E28:E32").Clear()

# ---------------------------------------------------------------------
# 3. Add the new course blocks.
# ---------------------------------------------------------------------

# GS-CC-6208, 10:00-11:00, N311 -> rows 10-14, columns B, D, F
$course1 = "GS-CC-6208 `n10:00-11:00 `nN311 "
foreach ($col in @("B","D","F")) {
    $rng = $ws.Range("$col" + "10:" + "$col" + "14")
    $rng.Merge()
    $rng.HorizontalAlignment = -4108
    $rng.WrapText = $true
}
$ws.Range("B10").Value = $course1
$ws.Range("D10").Value = $course1
$ws.Range("F10").Value = $course1

# GS-CC-6202, 11:00-12:30, M616 -> rows 14-20, column E (simple merge)
$course2 = "GS-CC-6202 `n11:00-12:30 `nM616 "
$rngE = $ws.Range("E14:E20")
$rngE.Merge()
$rngE.HorizontalAlignment = -4108
$rngE.WrapText = $true
$ws.Range("E14").Value = $course2

# Column C: this block was first merged across C14:C22, then narrowed to
# C14:C20, leaving the two bottom rows (C21:C22) formatted (centered /
# wrapped) but no longer part of the merged area.
$rngC = $ws.Range("C14:C22")
$rngC.Merge()
$rngC.HorizontalAlignment = -4108
$rngC.WrapText = $true
$ws.Range("C14").Value = $course2
$ws.Range("C21:C22").UnMerge()
$ws.Range("C14:C20").Merge()

# ---------------------------------------------------------------------
# 4. Reset row heights that Excel auto-grew because of the embedded
#    line breaks in the cell text, so the rows stay at their default
#    height (matches the original workbook's row styling).
# ---------------------------------------------------------------------
$ws.Rows("10:22").AutoFit()

Write-Host "Calendar updated."
